# WAT new script implementation
# - WAT23 (row 24) JIRA ID gains an additional linked ticket: WAT-142 -> WAT-142||WAT-153
# - WAT24 (row 25) JIRA ID gains an additional linked ticket: WAT-548 -> WAT-548||WAT-170
# - A brand new test case WAT26 is appended as row 31 (publication card / "more than 1" variant
#   of the WAT23 cart scenario)

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test Cases")

# --- Update existing JIRA ID references (WAT23 / WAT24 rows) ---
$ws.Range("B24").Value = "WAT-142||WAT-153"
$ws.Range("B25").Value = "WAT-548||WAT-170"

# --- Append the new WAT26 test case row ---
$newRow = 31
$ws.Range("A" + $newRow).Value = "WAT26"
$ws.Range("B" + $newRow).Value = "WAT-178"
$ws.Range("C" + $newRow).Value = "Verify that system must display following publication details in card if publication details is morethan 1, Publication count with morethan 1, Years, Top Journals, Recent publications link"
$ws.Range("D" + $newRow).Value = "Y"

# Formatting to match the rest of the table: thin borders around the row,
# wrapped description text, and a taller row (consistent with the other
# multi-line description rows such as 24 and 30).
$rowRange = $ws.Range("A" + $newRow + ":E" + $newRow)
$rowRange.Borders.LineStyle = 1
$ws.Range("C" + $newRow).WrapText = $true
$ws.Rows.Item($newRow).RowHeight = 30

# --- Update the view state to reflect scrolling/selection after the edit ---
$ws.Range("B35").Select() | Out-Null
